# Leetcode 75 workbook update: "Top K Frequent Elements Complete"
#
# - Row 33 (Pacific Atlantic Water Flow): Topic Heap -> Graph
# - Row 35 (Top K Frequent Elements): Topic Intervals -> Heap, mark Completed,
#   fill in Solution Notes + Runtime
# - Row 38 (Non-overlapping Intervals): Topic Linked List -> Intervals
# - Update the active selection to reflect where the author ended up (F35)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blind 75")

# --- Topic corrections -------------------------------------------------
$ws.Range("C33").Value = "Graph"
$ws.Range("C35").Value = "Heap"
$ws.Range("C38").Value = "Intervals"

# --- Top K Frequent Elements: mark complete with notes ------------------
$ws.Range("E35").Value = "X"
$ws.Range("G35").Value = "O(Nlog(k))"
$ws.Range("F35").Value = "count occurance of each num in dictionary. Use priority qeue to pop top k elements by value."

# --- Window / selection state -------------------------------------------
$win = $excel.ActiveWindow
[void]($win.ScrollRow = 22)
[void]($win.ScrollColumn = 2)
[void]$ws.Range("F35").Select()
